# Refresh the cryptos list (prices + 1h volume change %) with the latest
# snapshot. A couple of rows (InternetComputer / WEMIXToken) also swapped
# rank order, so Coin/Link/Price/Volume are rewritten together for those.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.305.00"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.584.58"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'209.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.807.59"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "1.582.72"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "26.308.83"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'7.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'206.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("D24").Value = "'8.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "'144.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'7.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.04%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "1.283.78"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "'0.611"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "'0.0167"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'0.767"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "'62.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "1.719.92"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "'88.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").Value = "0.0₇0968"
$ws.Range("E50").Value = "  -7.00%  "
$ws.Range("E51").Value = "  -0.11%  "
